$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.71609999999998
$ws.Range("A7").Value = -19.51419999999999
$ws.Range("C7").Value = -12.71690000000001
$ws.Range("C15").Value = -14.37619999999999
$ws.Range("A16").Value = -22.09940000000001
$ws.Range("E16").Value = 16.4022
$ws.Range("E19").Value = 16.39129999999999
$ws.Range("C21").Value = -12.2595
$ws.Range("C22").Value = -12.45070000000001
$ws.Range("C23").Value = -12.4798
$ws.Range("A28").Value = -21.91189999999998
$ws.Range("A29").Value = -21.35589999999998
$ws.Range("A32").Value = -21.2521
$ws.Range("C34").Value = -11.48780000000001
$ws.Range("E36").Value = 15.92490000000001
$ws.Range("A40").Value = -20.07539999999999
$ws.Range("C43").Value = -13.25649999999998
$ws.Range("C45").Value = -13.67909999999999
$ws.Range("E46").Value = 17.17279999999998
$ws.Range("C50").Value = -14.17239999999999
$ws.Range("E50").Value = 16.46079999999999
$ws.Range("C51").Value = -11.8827
$ws.Range("A52").Value = -22.28430000000001
$ws.Range("A57").Value = -22.36260000000001
$ws.Range("A66").Value = -21.99609999999999
$ws.Range("C66").Value = -12.306
$ws.Range("C67").Value = -10.80490000000001
$ws.Range("C79").Value = -11.45060000000001
$ws.Range("C84").Value = -13.5545
$ws.Range("C92").Value = -11.28360000000001
$ws.Range("E95").Value = 18.23010000000001
$ws.Range("C97").Value = -11.10760000000001
$ws.Range("E97").Value = 16.7204
$ws.Range("A100").Value = -22.03319999999999
